$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct urbansim_path (and corresponding urbansim_runid) for the
# 2035_TM152_NGF_NoProject_01 run (row 8) in the model log.
$ws.Range("F8").Value = "`"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.25 - FINAL VERSION`""
$ws.Range("G8").Value = "run182"

$ws.Range("F8").Select()
